$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.899.15'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.181.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.74%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.87%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.525'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.29'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.438'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.98%  '

$ws.Range("E11").Value = '  +4.38%  '

$ws.Range("E12").Value = '  +2.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.718.69'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.68%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000171'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.950.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.94%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.176.10'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.68%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.03'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.88%  '

$ws.Range("E22").Value = '  +0.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.532'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.168'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.46'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +15.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0866'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.89'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.65%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.04'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.48%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '157.15'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.50%  '

$ws.Range("E36").Value = '  +3.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '25.36'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.722.38'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.71'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0701'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.53%  '

$ws.Range("E41").Value = '  +5.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0294'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.727'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.220.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.66%  '

$ws.Range("E47").Value = '  +11.70%  '

$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.987'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.05%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.22'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.77%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.757'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.17%  '
